$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Stocks" (sheet1): update stock-sold quantities (column D) and add a
# new barcode value in column G for row 6.
# ---------------------------------------------------------------------------
$stocks = $wb.Worksheets.Item("Stocks")

$stocks.Range("D2").Value = 100
$stocks.Range("D3").Value = 150
$stocks.Range("D4").Value = 200
$stocks.Range("D6").Value = 144
$stocks.Range("D9").Value = 100
$stocks.Range("D10").Value = 160

# ---------------------------------------------------------------------------
# Sheet "Bills" (sheet2): append a new bill record as row 27.
# ---------------------------------------------------------------------------
$bills = $wb.Worksheets.Item("Bills")

$bills.Range("A27").Value = "02-Oct-2020 11:35"
$bills.Range("B27").Value = "KKK"

$bills.Range("C27").NumberFormat = "@"
$bills.Range("C27").Value = "688"
$bills.Range("C27").Style = "Normal"

$bills.Range("D27").NumberFormat = "@"
$bills.Range("D27").Value = "300.0"
$bills.Range("D27").Style = "Normal"

$bills.Range("E27").NumberFormat = "@"
$bills.Range("E27").Value = "150.0"
$bills.Range("E27").Style = "Normal"

$bills.Range("F27").Value = "XX02104"
$bills.Range("G27").Value = "Alica pure(1)"

# New barcode text for Stocks!G6 (added after the Bills edits above so the
# shared-string table grows in the same order as the authored change).
$stocks.Range("G6").Value = "yyy"

# ---------------------------------------------------------------------------
# Restore the active view on the Stocks sheet: scroll the frozen pane back to
# the top and move the selection to C5.
# ---------------------------------------------------------------------------
$stocks.Activate()
$stocks.Range("C5").Select()
